# Foodies_Project_Review_Template.xlsx — "last SIQ Reviewing" update.
#
# The four oldest reviews (rows 2-5 / RVW-001..RVW-004) had their Status
# column still sitting on "Pending". This SIQ review is the last one
# before hand-off to the customer, so every reviewed row's Status gets
# flipped to "Approved" (rows 6-10 already say "Approved" and are left
# untouched). With no cell left referencing the string "Pending", it
# naturally drops out of the workbook's shared string table.
#
# Finally, the view is reset to rest on the freshly-updated J2 cell
# instead of the previous J10 selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = "Approved"
$ws.Range("J3").Value = "Approved"
$ws.Range("J4").Value = "Approved"
$ws.Range("J5").Value = "Approved"

[void]$ws.Range("J2").Select()
